$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Update "Тип счетчика" (meter type) values for the example rows to showcase
# the new available meter types: ЭЛ (electricity), ТЕПЛО (heat), ГАЗ (gas)
$ws.Range("D6").Value = "ЭЛ"
$ws.Range("D7").Value = "ЭЛ"
$ws.Range("D8").Value = "ТЕПЛО"
$ws.Range("D9").Value = "ТЕПЛО"
$ws.Range("D10").Value = "ГАЗ"

# Move the active selection as in the authored workbook
$ws.Range("D16").Select() | Out-Null
